$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9272034764289856
$ws.Range("B1").Value = 1.671199917793274
$ws.Range("C1").Value = 4.467481136322021
$ws.Range("D1").Value = 2.253851652145386
$ws.Range("E1").Value = 0.9190675020217896
